$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.274017668686611
$ws.Range("D2").Value = 0.2710941420829158
$ws.Range("E2").Value = 0.1598654156039032
$ws.Range("F2").Value = 5.80377835155889
$ws.Range("G2").Value = 0.002665353917939985
$ws.Range("J2").Value = 0.1588548528935405
$ws.Range("L2").Value = 1.770097524809216
$ws.Range("N2").Value = 2.177618771485403
$ws.Range("B3").Value = 3.181425018440905
$ws.Range("D3").Value = 0.242474539260229
$ws.Range("E3").Value = 0.1389773348360919
$ws.Range("F3").Value = 5.727974780105598
$ws.Range("G3").Value = 0.00267629775453989
$ws.Range("J3").Value = 0.1388565524998455
$ws.Range("L3").Value = 1.68881301134553
$ws.Range("N3").Value = 2.207974208628858
$ws.Range("B4").Value = 3.126601116329141
$ws.Range("D4").Value = 0.2250934461947907
$ws.Range("E4").Value = 0.1261987174310235
$ws.Range("F4").Value = 5.685594491156053
$ws.Range("G4").Value = 0.002683356131724709
$ws.Range("J4").Value = 0.1265694494395433
$ws.Range("L4").Value = 1.640109468778036
$ws.Range("N4").Value = 2.227450041071972
$ws.Range("B5").Value = 3.104765837960485
$ws.Range("D5").Value = 0.218056357906903
$ws.Range("E5").Value = 0.1210017800359893
$ws.Range("F5").Value = 5.669359958579889
$ws.Range("G5").Value = 0.002686318072468217
$ws.Range("J5").Value = 0.1215593185319364
$ws.Range("L5").Value = 1.620560928572331
$ws.Range("N5").Value = 2.235596785296302
$ws.Range("B6").Value = 3.10117054362496
$ws.Range("D6").Value = 0.2168905576006352
$ws.Range("E6").Value = 0.1201394263107431
$ws.Range("F6").Value = 5.666726464126214
$ws.Range("G6").Value = 0.002686815081825156
$ws.Range("J6").Value = 0.1207271776583525
$ws.Range("L6").Value = 1.617332808472526
$ws.Range("N6").Value = 2.236962229036664
$ws.Range("B7").Value = 3.126304595098588
$ws.Range("D7").Value = 0.2249983589254612
$ws.Range("E7").Value = 0.1261285890950745
$ws.Range("F7").Value = 5.685371367268488
$ws.Range("G7").Value = 0.002683395730320092
$ws.Range("J7").Value = 0.1265018946744334
$ws.Range("L7").Value = 1.639844627472939
$ws.Range("N7").Value = 2.227559060400999
$ws.Range("B8").Value = 3.241668170978073
$ws.Range("D8").Value = 0.2611849629434175
$ws.Range("E8").Value = 0.1526527051487747
$ws.Range("F8").Value = 5.776770348156504
$ws.Range("G8").Value = 0.00266905727114173
$ws.Range("J8").Value = 0.151960493949332
$ws.Range("L8").Value = 1.741817608507972
$ws.Range("N8").Value = 2.187911265799141
$ws.Range("B9").Value = 3.484197717918562
$ws.Range("D9").Value = 0.3337745525760738
$ws.Range("E9").Value = 0.2050995078894289
$ws.Range("F9").Value = 5.989567191521104
$ws.Range("G9").Value = 0.002643609789433698
$ws.Range("J9").Value = 0.2018693338227138
$ws.Range("L9").Value = 1.951574718333006
$ws.Range("N9").Value = 2.116824163447774
$ws.Range("B10").Value = 3.672640947273976
$ws.Range("D10").Value = 0.3882564389307106
$ws.Range("E10").Value = 0.2439873438326288
$ws.Range("F10").Value = 6.167141899548511
$ws.Range("G10").Value = 0.002626515765267623
$ws.Range("J10").Value = 0.2385998833442784
$ws.Range("L10").Value = 2.111985074940037
$ws.Range("N10").Value = 2.068680068767478
$ws.Range("B11").Value = 3.760669081248864
$ws.Range("D11").Value = 0.4133269312012828
$ws.Range("E11").Value = 0.2617760934801083
$ws.Range("F11").Value = 6.252709816980541
$ws.Range("G11").Value = 0.002619081621218823
$ws.Range("J11").Value = 0.2553394092936685
$ws.Range("L11").Value = 2.186404504684333
$ws.Range("N11").Value = 2.04766984203016
$ws.Range("B12").Value = 3.79434017648839
$ws.Range("D12").Value = 0.4228645494701766
$ws.Range("E12").Value = 0.2685280520861681
$ws.Range("F12").Value = 6.285814985667287
$ws.Range("G12").Value = 0.002616315249320049
$ws.Range("J12").Value = 0.2616839405433495
$ws.Range("L12").Value = 2.214799529581398
$ws.Range("N12").Value = 2.039842427182702
$ws.Range("B13").Value = 3.787073438647553
$ws.Range("D13").Value = 0.4208084521574165
$ws.Range("E13").Value = 0.2670731709945642
$ws.Range("F13").Value = 6.278653727566763
$ws.Range("G13").Value = 0.002616908873947689
$ws.Range("J13").Value = 0.2603172614940092
$ws.Range("L13").Value = 2.208674540845095
$ws.Range("N13").Value = 2.041522466986358
$ws.Range("B14").Value = 3.763432441962038
$ws.Range("D14").Value = 0.4141106994173924
$ws.Range("E14").Value = 0.2623312563770526
$ws.Range("F14").Value = 6.255419239788949
$ws.Range("G14").Value = 0.002618853054888628
$ws.Range("J14").Value = 0.2558612575344625
$ws.Range("L14").Value = 2.188736257696121
$ws.Range("N14").Value = 2.047023292944257
$ws.Range("B15").Value = 3.748995667118834
$ws.Range("D15").Value = 0.4100139421557003
$ws.Range("E15").Value = 0.2594287950607992
$ws.Range("F15").Value = 6.241279345237103
$ws.Range("G15").Value = 0.00262005026176305
$ws.Range("J15").Value = 0.2531325977093104
$ws.Range("L15").Value = 2.176551528086122
$ws.Range("N15").Value = 2.050409487003428
$ws.Range("B16").Value = 3.666934927309796
$ws.Range("D16").Value = 0.3866240226110165
$ws.Range("E16").Value = 0.2428269231981943
$ws.Range("F16").Value = 6.161647366921414
$ws.Range("G16").Value = 0.002627008451681686
$ws.Range("J16").Value = 0.2375066365852945
$ws.Range("L16").Value = 2.107151217023784
$ws.Range("N16").Value = 2.070071112674455
$ws.Range("B17").Value = 3.617187294196128
$ws.Range("D17").Value = 0.3723503757632045
$ws.Range("E17").Value = 0.2326685476991912
$ws.Range("F17").Value = 6.114031889289635
$ws.Range("G17").Value = 0.00263136438687139
$ws.Range("J17").Value = 0.2279293224287358
$ws.Range("L17").Value = 2.064951286352652
$ws.Range("N17").Value = 2.082361456855502
$ws.Range("B18").Value = 3.588790334913767
$ws.Range("D18").Value = 0.3641673387896844
$ws.Range("E18").Value = 0.2268348904730573
$ws.Range("F18").Value = 6.087094566740262
$ws.Range("G18").Value = 0.002633902022018136
$ws.Range("J18").Value = 0.2224235101354282
$ws.Range("L18").Value = 2.040814932714397
$ws.Range("N18").Value = 2.089514342728666
$ws.Range("B19").Value = 3.579212646463361
$ws.Range("D19").Value = 0.3614012203136951
$ws.Range("E19").Value = 0.2248612421178677
$ws.Range("F19").Value = 6.078050895589996
$ws.Range("G19").Value = 0.00263476676708309
$ws.Range("J19").Value = 0.2205597815914047
$ws.Range("L19").Value = 2.032665951082606
$ws.Range("N19").Value = 2.091950562263751
$ws.Range("B20").Value = 3.622460562788149
$ws.Range("D20").Value = 0.3738670375607001
$ws.Range("E20").Value = 0.2337489635499423
$ws.Range("F20").Value = 6.119053981945939
$ws.Range("G20").Value = 0.002630897358887364
$ws.Range("J20").Value = 0.22894854469682
$ws.Range("L20").Value = 2.069429431440199
$ws.Range("N20").Value = 2.081044450889241
$ws.Range("B21").Value = 3.770367192483775
$ws.Range("D21").Value = 0.4160767759514101
$ws.Range("E21").Value = 0.263723631470981
$ws.Range("F21").Value = 6.262224594750762
$ws.Range("G21").Value = 0.002618280681341899
$ws.Range("J21").Value = 0.2571699309281712
$ws.Range("L21").Value = 2.194586759706851
$ws.Range("N21").Value = 2.045404069090512
$ws.Range("B22").Value = 3.868998214399312
$ws.Range("D22").Value = 0.443920900098874
$ws.Range("E22").Value = 0.2834063536658817
$ws.Range("F22").Value = 6.359896063876477
$ws.Range("G22").Value = 0.002610319073694853
$ws.Range("J22").Value = 0.2756476901399481
$ws.Range("L22").Value = 2.277634504650507
$ws.Range("N22").Value = 2.022861493371128
$ws.Range("B23").Value = 3.816175306937453
$ws.Range("D23").Value = 0.4290354552786653
$ws.Range("E23").Value = 0.272892310927304
$ws.Range("F23").Value = 6.30738700070782
$ws.Range("G23").Value = 0.002614542470560797
$ws.Range("J23").Value = 0.2657822836360424
$ws.Range("L23").Value = 2.233194003009828
$ws.Range("N23").Value = 2.034824013872999
$ws.Range("B24").Value = 3.620075884238474
$ws.Range("D24").Value = 0.3731812831045147
$ws.Range("E24").Value = 0.233260487329261
$ws.Range("F24").Value = 6.116782132344582
$ws.Range("G24").Value = 0.00263110839837507
$ws.Range("J24").Value = 0.2284877534268475
$ws.Range("L24").Value = 2.067404473293948
$ws.Range("N24").Value = 2.081639598336935
$ws.Range("B25").Value = 3.416809514448516
$ws.Range("D25").Value = 0.313946792668105
$ws.Range("E25").Value = 0.1908556167558118
$ws.Range("F25").Value = 5.928328319621272
$ws.Range("G25").Value = 0.00265021078121618
$ws.Range("J25").Value = 0.1883620061578455
$ws.Range("L25").Value = 1.893747989492681
$ws.Range("N25").Value = 2.041522466986358
